# Applies the diff between the original "current_rate" sheet values and the
# updated ones. All edits are plain value overwrites on the single sheet
# (no formulas are involved), plus a handful of cell deletions (S20:S23,
# V23:W23) and two new cells (V18, W18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("P2").Value2 = 0
$ws.Range("S2").Value2 = 72056
$ws.Range("T2").Value2 = 5516.6335
$ws.Range("U2").Value2 = 45000

# --- Row 3 ---
$ws.Range("P3").Value2 = 0
$ws.Range("S3").Value2 = 69308
$ws.Range("T3").Value2 = 5387.312
$ws.Range("U3").Value2 = 45000

# --- Row 4 ---
$ws.Range("B4").Value2 = 12500
$ws.Range("I4").Value2 = 10000
$ws.Range("P4").Value2 = 0
$ws.Range("S4").Value2 = 65612
$ws.Range("T4").Value2 = 5339.824000000001
$ws.Range("U4").Value2 = 22500

# --- Row 5 ---
$ws.Range("B5").Value2 = 12500
$ws.Range("I5").Value2 = 10000
$ws.Range("S5").Value2 = 64269
$ws.Range("T5").Value2 = 5306.819000000001
$ws.Range("U5").Value2 = 22500

# --- Row 6 ---
$ws.Range("B6").Value2 = 12500
$ws.Range("S6").Value2 = 65136
$ws.Range("T6").Value2 = 5338.83
$ws.Range("U6").Value2 = 32500

# --- Row 7 ---
$ws.Range("S7").Value2 = 71065
$ws.Range("T7").Value2 = 5422.714499999999

# --- Row 8 ---
$ws.Range("S8").Value2 = 70049
$ws.Range("T8").Value2 = 5875.107

# --- Row 9 ---
$ws.Range("S9").Value2 = 80484
$ws.Range("T9").Value2 = 7214.4275

# --- Row 10 ---
$ws.Range("B10").Value2 = 25000
$ws.Range("I10").Value2 = 20000
$ws.Range("P10").Value2 = 20000
$ws.Range("S10").Value2 = 94725
$ws.Range("T10").Value2 = 8467.018
$ws.Range("U10").Value2 = 65000

# --- Row 11 ---
$ws.Range("B11").Value2 = 25000
$ws.Range("I11").Value2 = 20000
$ws.Range("P11").Value2 = 20000
$ws.Range("S11").Value2 = 104123
$ws.Range("T11").Value2 = 13640.6165
$ws.Range("U11").Value2 = 65000

# --- Row 12 ---
$ws.Range("P12").Value2 = 20000
$ws.Range("S12").Value2 = 110312
$ws.Range("T12").Value2 = 15763.02
$ws.Range("U12").Value2 = 65000

# --- Row 13 ---
$ws.Range("P13").Value2 = 20000
$ws.Range("S13").Value2 = 111321
$ws.Range("T13").Value2 = 15072.0395
$ws.Range("U13").Value2 = 65000

# --- Row 14 ---
$ws.Range("P14").Value2 = 20000
$ws.Range("S14").Value2 = 109259
$ws.Range("T14").Value2 = 15259.6675
$ws.Range("U14").Value2 = 65000

# --- Row 15 ---
$ws.Range("S15").Value2 = 113429
$ws.Range("T15").Value2 = 15815.744

# --- Row 16 ---
$ws.Range("S16").Value2 = 114542
$ws.Range("T16").Value2 = 15955.3275

# --- Row 17 ---
$ws.Range("S17").Value2 = 99426
$ws.Range("T17").Value2 = 16114.0735

# --- Row 18 ---
$ws.Range("S18").Value2 = 92477
$ws.Range("T18").Value2 = 16384.389
$ws.Range("V18").Value2 = 5260.194847222222
$ws.Range("W18").Value2 = 5.241048973938897

# --- Row 19 ---
$ws.Range("S19").Value2 = 1237
$ws.Range("T19").Value2 = 15932.637

# --- Row 20 --- (S20 removed entirely)
$ws.Range("S20").ClearContents()
$ws.Range("T20").Value2 = 15036.9065

# --- Row 21 --- (S21 removed entirely)
$ws.Range("S21").ClearContents()
$ws.Range("T21").Value2 = 13330.5025

# --- Row 22 --- (S22 removed entirely)
$ws.Range("S22").ClearContents()
$ws.Range("T22").Value2 = 11770.8115

# --- Row 23 --- (S23, V23, W23 removed entirely)
$ws.Range("S23").ClearContents()
$ws.Range("T23").Value2 = 9312.075499999999
$ws.Range("V23").ClearContents()
$ws.Range("W23").ClearContents()

# --- Row 24 ---
$ws.Range("T24").Value2 = 6773.7985

# --- Row 25 ---
$ws.Range("T25").Value2 = 5542.130999999999
